$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '31.296.05'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.17%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.004.45'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +7.25%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9986'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.12%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7961'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +68.72%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '258.75'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.05%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9987'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.12%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3612'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +25.68%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '28.50'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +30.75%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07094'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +9.43%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8593'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +18.50%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08184'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.02%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.006.19'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.29%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '101.52'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.91%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.614'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +8.74%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '275.73'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.59%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '15.02'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +14.77%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '31.311.85'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.26%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.914'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +12.22%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007996'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.93%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.269.50'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.45%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9992'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.06%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9974'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.22%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.198'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +14.54%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.15'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +12.43%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1531'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +58.62%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.26%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.05'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.59%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.399'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +26.99%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.621'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +9.23%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.625'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +9.37%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.365'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.57%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.425'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.07%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05247'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +9.27%  '

# Row 35
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7807'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +13.34%  '

# Row 36
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.224'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.92%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.807'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.43%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02013'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.24%  '

# Row 39
$ws.Range("E39").Value = '  +3.55%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.733'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.43%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '80.76'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.13%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4777'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +13.31%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.160'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +10.71%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '107.49'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.62%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8589'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.17%  '

# Row 46
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.821'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +11.83%  '

# Row 47
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9993'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.02%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.955'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.37%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4374'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +12.27%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.96'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.65%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1196'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +14.53%  '
